# Add a "UnitEffect" column (visual effect on trap) to the SpellTrap table,
# plus a new trap row (54000007 / 冰冻陷阱), per the commit:
# "add the visual effect on trap"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- add the new "UnitEffect" column to the table -------------------------
$newCol = $lo.ListColumns.Add()

# carry over the header formatting from column E (Comment) to column F
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("E3").Copy()
$ws.Range("F3").PasteSpecial(-4122)

$ws.Range("F1").Value = "特效"
$ws.Range("F2").Value = "string"
$ws.Range("F3").Value = "UnitEffect"

# --- populate the new column for the existing rows -------------------------
$ws.Range("F4").Value = "bluewing"
$ws.Range("F5").Value = "silent"
$ws.Range("F6").Value = "silent"
$ws.Range("F7").Value = "pinkball"
$ws.Range("F8").Value = "icesharp"
$ws.Range("F9").Value = "firehit"

# --- add the new trap row (54000007 / Frost Trap) --------------------------
$newRow = $lo.ListRows.Add()
$ws.Rows.Item(10).RowHeight = 24

$ws.Range("A10").Value = 54000007
$ws.Range("B10").Value = "冰冻陷阱"
$ws.Range("C10").Value = "return false;"
$ws.Range("D10").Value = "m.Return((int)t.Help);return true;"
$ws.Range("E10").Value = "把目标移动回手牌"
$ws.Range("F10").Value = "iceball"

# --- match the selection left by the author in the source workbook --------
$ws.Range("F4").Select()
